$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.548.84'
$ws.Range('E2').Value = '  -2.30%  '
$ws.Range('D3').Value = '3.481.58'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'581.95"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.78%  '
$ws.Range('D6').Value = "'131.61"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.12%  '
$ws.Range('D7').Value = '3.483.94'
$ws.Range('E7').Value = '  -1.26%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'0.484"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.57%  '
$ws.Range('D10').Value = "'0.122"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.97%  '
$ws.Range('D11').Value = "'7.12"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('E12').Value = '  -1.00%  '
$ws.Range('D13').Value = '4.098.80'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('D14').Value = "'27.38"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('E15').Value = '  +1.65%  '
$ws.Range('D16').Value = '3.504.69'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('E17').Value = '  -3.51%  '
$ws.Range('D18').Value = '63.697.65'
$ws.Range('E18').Value = '  -0.44%  '
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').Value = "'14.28"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('D22').Value = "'382.87"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').Value = '3.626.95'
$ws.Range('E24').Value = '  -1.08%  '
$ws.Range('D25').Value = "'72.70"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.14%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = "'0.0000112"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('D28').Value = "'1.57"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.18%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = "'1.00"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').Value = "'7.45"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.99%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = "'8.28"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.08%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = "'2.24"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('D33').Value = '3.490.44'
$ws.Range('E33').Value = '  -1.07%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').Value = "'23.44"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.94%  '
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('D37').Value = "'5.31"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.31%  '
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('E39').Value = '  +1.30%  '
$ws.Range('D40').Value = "'161.27"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.42%  '
$ws.Range('E41').Value = '  -2.77%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').Value = "'0.808"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.67%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = "'26.28"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.95%  '
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').Value = "'1.21"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.48%  '
$ws.Range('D46').Value = "'41.27"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.32%  '
$ws.Range('D47').Value = "'4.37"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('D48').Value = "'1.62"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = "'6.81"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.42%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.412.82'
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('D51').Value = "'0.889"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.95%  '
